# Update "想去人数" (F) counts and one time-range (E24) across both the
# "展览" and "全部类型" worksheets, which carry duplicate data.

$wb = $excel.ActiveWorkbook

# Map of row -> new F-column value for both affected sheets.
$fUpdates = @{
    3  = 3025
    5  = 156
    7  = 1655
    11 = 1356
    13 = 489
    14 = 345
    15 = 22
    18 = 121
    21 = 3146
    23 = 118
    24 = 193
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $fUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $fUpdates[$row]
    }

    # Row 24's specific time range (column E) changes end time 17:00 -> 18:00
    $ws.Cells.Item(24, 5).Value = "2024.05.01 10:00-05.02 18:00"
}
